$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update station name column (D) with nicely formatted "City, ST" labels.
# Doing this in row order (2..47) so the shared-string table appends the
# new strings in the same order as the target workbook.
$ws.Range("D2").Value = "Rapid City NWS, SD"
$ws.Range("D3").Value = "Ellsworth AFB, SD"
$ws.Range("D4").Value = "Rapid City Airport, SD"
$ws.Range("D5").Value = "Rapid City, NEXRAD, SD"
$ws.Range("D6").Value = "Custer Airport, SD"
$ws.Range("D7").Value = "Clyde Ice Field, SD"
$ws.Range("D8").Value = "Belle Fourche Airport, SD"
$ws.Range("D9").Value = "Hulett Airport, WY"
$ws.Range("D10").Value = "Philip Airport, SD"
$ws.Range("D11").Value = "Pine Ridge Airport, SD"
$ws.Range("D12").Value = "Chadron Airport, NE"
$ws.Range("D13").Value = "Moorcroft, WY"
$ws.Range("D14").Value = "Faith Airport, SD"
$ws.Range("D15").Value = "Gordon Airport, NE"
$ws.Range("D16").Value = "Buffalo, SD"
$ws.Range("D17").Value = "Gillette  Airport, WY"
$ws.Range("D18").Value = "Hettinger Airport, ND"
$ws.Range("D19").Value = "Lemmon, SD"
$ws.Range("D20").Value = "Lemmon Airport, SD"
$ws.Range("D21").Value = "Converse, CO Airport, WY"
$ws.Range("D22").Value = "Alliance Airport, NE"
$ws.Range("D23").Value = "Douglas, WY"
$ws.Range("D24").Value = "Bowman Airport, ND"
$ws.Range("D25").Value = "Sibley Peak, WY"
$ws.Range("D26").Value = "Pierre Airport, SD"
$ws.Range("D27").Value = "Torrington Airport, WY"
$ws.Range("D28").Value = "Mission, SD"
$ws.Range("D29").Value = "Scottsbluff Airport, NE"
$ws.Range("D30").Value = "Miller Field Airport, NE"
$ws.Range("D31").Value = "Cheyenne NWS, WY"
$ws.Range("D32").Value = "Bismark NWS, ND"
$ws.Range("D33").Value = "North Platte NWS, NE"
$ws.Range("D34").Value = "Aberdeen NWS, SD"
$ws.Range("D35").Value = "Riverton NWS, WY"
$ws.Range("D36").Value = "Billings NWS, MT"
$ws.Range("D37").Value = "Denver Intl Airport, CO"
$ws.Range("D38").Value = "Boulder NWS, CO"
$ws.Range("D39").Value = "Sioux Falls NWS, SD"
$ws.Range("D40").Value = "Glasgow NWS, MT"
$ws.Range("D41").Value = "Goodland NWS KS"
$ws.Range("D42").Value = "Hastings NWS, NE"
$ws.Range("D43").Value = "Munich School, ND"
$ws.Range("D44").Value = "Grand Forks NWS, ND"
$ws.Range("D45").Value = "Omaha-Valley NWS, NE"
$ws.Range("D46").Value = "Grand Junction NWS, CO"
$ws.Range("D47").Value = "Topeka NWS, KS"

# Update the active selection/view to match the reviewed state: D2:D47
# selected with D2 as the active cell, scrolled near the bottom of the list.
$ws.Range("D2:D47").Select()
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1
